# Add the missing "Time out" punch for the Oct 13 2024 entry (row 15) and
# roll it into the Week-3 subtotal (rows 10-15) on the log sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15 was clocked in (B15) but never clocked out (C15 was blank).
# Fill in the clock-out time (13:05 -> serial fraction of a day).
$ws.Range("C15").Value = 0.54513888888888884

# Week-3 totals previously summed rows 10:14; extend them to include the
# newly completed row 15.
$ws.Range("M4").Formula = "=SUM(D10:D15)"
$ws.Range("N4").Formula = "=SUM(G10:G15)"

# Leave the selection where the user ended up after entering the time.
$ws.Range("C16").Select() | Out-Null
